$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.136730802768966
$ws.Range("C2").Value = 0.1710915350730318
$ws.Range("D2").Value = 0.05164015962787971
$ws.Range("E2").Value = 0.1237518460948799
$ws.Range("F2").Value = 3.29784889486362
$ws.Range("J2").Value = 0.2261291256545519
$ws.Range("K2").Value = 1.229727015956257
$ws.Range("N2").Value = 3.111090195120724
$ws.Range("B3").Value = 1.09089379475634
$ws.Range("C3").Value = 0.1629003829944224
$ws.Range("D3").Value = 0.05030122963908212
$ws.Range("E3").Value = 0.1211336060323056
$ws.Range("F3").Value = 3.273361124162491
$ws.Range("J3").Value = 0.221685578417059
$ws.Range("K3").Value = 1.178190574632794
$ws.Range("N3").Value = 3.11691075222646
$ws.Range("B4").Value = 1.063452048856959
$ws.Range("C4").Value = 0.1579826333690448
$ws.Range("D4").Value = 0.04951307680558159
$ws.Range("E4").Value = 0.1195983299882002
$ws.Range("F4").Value = 3.259770387011045
$ws.Range("J4").Value = 0.2190838162326401
$ws.Range("K4").Value = 1.147316752513689
$ws.Range("N4").Value = 3.121222168917299
$ws.Range("B5").Value = 1.052445612452317
$ws.Range("C5").Value = 0.1560065759215234
$ws.Range("D5").Value = 0.04920046178716575
$ws.Range("E5").Value = 0.1189908644581692
$ws.Range("F5").Value = 3.254594835344449
$ws.Range("J5").Value = 0.2180553356912256
$ws.Range("K5").Value = 1.134928550289231
$ws.Range("N5").Value = 3.123164253034972
$ws.Range("B6").Value = 1.050628641023849
$ws.Range("C6").Value = 0.1556801389632767
$ws.Range("D6").Value = 0.04914907027813342
$ws.Range("E6").Value = 0.1188910923313884
$ws.Range("F6").Value = 3.253757334788133
$ws.Range("J6").Value = 0.2178864732844019
$ws.Range("K6").Value = 1.132883151052198
$ws.Range("N6").Value = 3.123497907706607
$ws.Range("B7").Value = 1.063302898799009
$ws.Range("C7").Value = 0.1579558704798245
$ws.Range("D7").Value = 0.04950882606134144
$ws.Range("E7").Value = 0.1195900639395333
$ws.Range("F7").Value = 3.259699119536265
$ws.Range("J7").Value = 0.2190698172905243
$ws.Range("K7").Value = 1.147148899472484
$ws.Range("N7").Value = 3.121247611338845
$ws.Range("B8").Value = 1.120780216510809
$ws.Range("C8").Value = 0.1682439829778559
$ws.Range("D8").Value = 0.05117146351103941
$ws.Range("E8").Value = 0.1228340523967688
$ws.Range("F8").Value = 3.289105292169097
$ws.Range("J8").Value = 0.2245706690116833
$ws.Range("K8").Value = 1.21179720496869
$ws.Range("N8").Value = 3.112943843857693
$ws.Range("B9").Value = 1.239088032939776
$ws.Range("C9").Value = 0.1893111294020287
$ws.Range("D9").Value = 0.05470047065605144
$ws.Range("E9").Value = 0.1297707430806412
$ws.Range("F9").Value = 3.358265995471697
$ws.Range("J9").Value = 0.2363665225199156
$ws.Range("K9").Value = 1.344708142913873
$ws.Range("N9").Value = 3.102528086686831
$ws.Range("B10").Value = 1.329461333775271
$ws.Range("C10").Value = 0.2053436850826245
$ws.Range("D10").Value = 0.05745622924126792
$ws.Range("E10").Value = 0.1352203846029738
$ws.Range("F10").Value = 3.416139065275388
$ws.Range("J10").Value = 0.2456550031435825
$ws.Range("K10").Value = 1.446149712523265
$ws.Range("N10").Value = 3.098476087102242
$ws.Range("B11").Value = 1.371333998530304
$ws.Range("C11").Value = 0.2127601558477181
$ws.Range("D11").Value = 0.05874516512757566
$ws.Range("E11").Value = 0.1377768884714996
$ws.Range("F11").Value = 3.444012401755174
$ws.Range("J11").Value = 0.2500172867448072
$ws.Range("K11").Value = 1.493133486145524
$ws.Range("N11").Value = 3.097419738053688
$ws.Range("B12").Value = 1.387300156466154
$ws.Range("C12").Value = 0.2155864551942841
$ws.Range("D12").Value = 0.05923831415145742
$ws.Range("E12").Value = 0.1387561407712283
$ws.Range("F12").Value = 3.454790530479556
$ws.Range("J12").Value = 0.251688966708997
$ws.Range("K12").Value = 1.511046184117959
$ws.Range("N12").Value = 3.097133297988066
$ws.Range("B13").Value = 1.383856664859763
$ws.Range("C13").Value = 0.2149769654074021
$ws.Range("D13").Value = 0.059131881136949
$ws.Range("E13").Value = 0.1385447445302432
$ws.Range("F13").Value = 3.452459333764779
$ws.Range("J13").Value = 0.2513280594321401
$ws.Range("K13").Value = 1.507182975654075
$ws.Range("N13").Value = 3.097189929922166
$ws.Range("B14").Value = 1.372645339287601
$ws.Range("C14").Value = 0.2129923188492171
$ws.Range("D14").Value = 0.05878563557578786
$ws.Range("E14").Value = 0.1378572283441244
$ws.Range("F14").Value = 3.444894649531932
$ws.Range("J14").Value = 0.2501544199584487
$ws.Range("K14").Value = 1.494604747223889
$ws.Range("N14").Value = 3.097393893805787
$ws.Range("B15").Value = 1.365792401471083
$ws.Range("C15").Value = 0.2117789930781271
$ws.Range("D15").Value = 0.05857420806258062
$ws.Range("E15").Value = 0.1374375589100794
$ws.Range("F15").Value = 3.440290137495765
$ws.Range("J15").Value = 0.2494381102726919
$ws.Range("K15").Value = 1.48691598803353
$ws.Range("N15").Value = 3.097533631061509
$ws.Range("B16").Value = 1.326740211896038
$ws.Range("C16").Value = 0.204861491008046
$ws.Range("D16").Value = 0.05737270314558174
$ws.Range("E16").Value = 0.135054871342021
$ws.Range("F16").Value = 3.414348658304931
$ws.Range("J16").Value = 0.2453726802944232
$ws.Range("K16").Value = 1.443096110118546
$ws.Range("N16").Value = 3.098560992498633
$ws.Range("B17").Value = 1.302978251960553
$ws.Range("C17").Value = 0.2006494685102496
$ws.Range("D17").Value = 0.05664465043907541
$ws.Range("E17").Value = 0.133613019896579
$ws.Range("F17").Value = 3.398831015317938
$ws.Range("J17").Value = 0.2429137999594104
$ws.Range("K17").Value = 1.416428911069204
$ws.Range("N17").Value = 3.099393091108695
$ws.Range("B18").Value = 1.289382642731937
$ws.Range("C18").Value = 0.1982384110322926
$ws.Range("D18").Value = 0.05622922031599842
$ws.Range("E18").Value = 0.1327909920585313
$ws.Range("F18").Value = 3.390051206649815
$ws.Range("J18").Value = 0.2415123978559421
$ws.Range("K18").Value = 1.401169439868852
$ws.Range("N18").Value = 3.099945736105866
$ws.Range("B19").Value = 1.284791695547995
$ws.Range("C19").Value = 0.1974240544753627
$ws.Range("D19").Value = 0.05608913486275924
$ws.Range("E19").Value = 0.1325139183620152
$ws.Range("F19").Value = 3.387103487831581
$ws.Range("J19").Value = 0.241040116737878
$ws.Range("K19").Value = 1.396016362474057
$ws.Range("N19").Value = 3.100145556679706
$ws.Range("B20").Value = 1.305500335665215
$ws.Range("C20").Value = 0.2010966454101322
$ws.Range("D20").Value = 0.05672180877655109
$ws.Range("E20").Value = 0.1337657530927814
$ws.Range("F20").Value = 3.400467828184929
$ws.Range("J20").Value = 0.2431742181984617
$ws.Range("K20").Value = 1.419259521004562
$ws.Range("N20").Value = 3.09929684657331
$ws.Range("B21").Value = 1.375935391781582
$ws.Range("C21").Value = 0.2135747723570205
$ws.Range("D21").Value = 0.05888719922645436
$ws.Range("E21").Value = 0.1380588654390493
$ws.Range("F21").Value = 3.44711052007986
$ws.Range("J21").Value = 0.2504986088255379
$ws.Range("K21").Value = 1.498295990619937
$ws.Range("N21").Value = 3.097330898880173
$ws.Range("B22").Value = 1.422609481292113
$ws.Range("C22").Value = 0.2218340132270384
$ws.Range("D22").Value = 0.06033187264347362
$ws.Range("E22").Value = 0.1409297223354784
$ws.Range("F22").Value = 3.478894956549112
$ws.Range("J22").Value = 0.2554008451648997
$ws.Range("K22").Value = 1.550656211369471
$ws.Range("N22").Value = 3.096708246632772
$ws.Range("B23").Value = 1.397639881090868
$ws.Range("C23").Value = 0.2174163356423264
$ws.Range("D23").Value = 0.05955813464312598
$ws.Range("E23").Value = 0.1393915305878366
$ws.Range("F23").Value = 3.461811750428012
$ws.Range("J23").Value = 0.2527738468377692
$ws.Range("K23").Value = 1.522645860179352
$ws.Range("N23").Value = 3.096979836233729
$ws.Range("B24").Value = 1.304359897761174
$ws.Range("C24").Value = 0.2008944440688936
$ws.Range("D24").Value = 0.05668691571881368
$ws.Range("E24").Value = 0.1336966808883773
$ws.Range("F24").Value = 3.399727384403747
$ws.Range("J24").Value = 0.2430564449922343
$ws.Range("K24").Value = 1.417979578471488
$ws.Range("N24").Value = 3.099340127430366
$ws.Range("B25").Value = 1.206479383211473
$ws.Range("C25").Value = 0.1835153039756392
$ws.Range("D25").Value = 0.0537171053297314
$ws.Range("E25").Value = 0.1278323382274955
$ws.Range("F25").Value = 3.338320128606341
$ws.Range("J25").Value = 0.2330667112360203
$ws.Range("K25").Value = 1.308090071280191
$ws.Range("N25").Value = 3.104715246303087
